$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New cell D2: new shared string about "Falta en cuenta corriente - nueva cuota - ingreso de autorizacion"
$ws.Range("D2").Value = "Falta en cuenta corriente - nueva cuota - ingreso de autorizacion"

# New row 27 with two new strings
$ws.Range("A27").Value = "Validar fechas - permite cargar fechas q no existen"
$ws.Range("C27").Value = "Lucas - ver como agregar validacion de fecha al fwk de validación que hiciste"

# Update the view: scroll so row 19 is the top visible row, and select C28 (one row below the new last row)
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("C28").Select()
